# Add a new worksheet "Grades2" at the end of the workbook containing a
# single column of grade values (header "Grades" + 13 numeric grades),
# and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Insert the new sheet immediately after the current last sheet so it
# lands at the end of the tab strip (and becomes sheetId 5 / rId5).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "Grades2"

# Header
$ws.Range("A1").Value = "Grades"

# Data
$grades = @(8, 4.5, 4.4, 4.7, 5.8, 4.6, 3.9, 3.9, 8.2, 5.1, 7.1, 4.1, 5.8)
for ($i = 0; $i -lt $grades.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $grades[$i]
}

# Make the new sheet the active one, with B2 selected (matches the
# target's <selection activeCell="B2" sqref="B2"/> on the new sheet).
$ws.Activate() | Out-Null
$ws.Range("B2").Select() | Out-Null
